$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was added to the daily price log. It becomes the
# new row 424, pushing every existing row from the old 424 down through the
# old 482 down by one (to 425..483).
$ws.Rows.Item(424).Insert()

$ws.Cells.Item(424, 1).Value  = 10
$ws.Cells.Item(424, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(424, 3).Value  = "La Araucanía"
$ws.Cells.Item(424, 4).Value  = 45127
$ws.Cells.Item(424, 5).Value  = 9
$ws.Cells.Item(424, 6).Value  = 100112001
$ws.Cells.Item(424, 7).Value  = "Berenjena"
$ws.Cells.Item(424, 8).Value  = "Sin especificar"
$ws.Cells.Item(424, 9).Value  = "Primera"
$ws.Cells.Item(424, 10).Value = 150
$ws.Cells.Item(424, 11).Value = 10000
$ws.Cells.Item(424, 12).Value = 10000
$ws.Cells.Item(424, 13).Value = 10000
$ws.Cells.Item(424, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(424, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(424, 16).Value = 250
$ws.Cells.Item(424, 17).Value = 40
$ws.Cells.Item(424, 18).Value = "Hortaliza"
